$d = $word.ActiveDocument
$wmain = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark (tracks the location of the author's most
#    recent edit) from the "public CreateGameView( )" subtitle line into the
#    "Some paragraph about this class" filler paragraph, splitting the 4th
#    occurrence of that sentence into "Some paragraph ab" | "out this class".
# ---------------------------------------------------------------------------
$paraFiller = $d.Paragraphs.Item(13)
$fullText = $paraFiller.Range.Text
$marker = "Some paragraph about this class"
$searchFrom = 0
$occurrence = 0
$targetOccurrence = 4
$foundPos = -1
while ($true) {
    $idx = $fullText.IndexOf($marker, $searchFrom)
    if ($idx -lt 0) { break }
    $occurrence++
    if ($occurrence -eq $targetOccurrence) {
        $foundPos = $idx
        break
    }
    $searchFrom = $idx + 1
}
$splitOffset = $foundPos + ("Some paragraph ab").Length
$absOffset = $paraFiller.Range.Start + $splitOffset
$bmRange = $d.Range($absOffset, $absOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 2) "public CreateGameView( )" subtitle paragraph: mark the constructor
#    signature with gramStart/gramEnd (matching the Heading4 write-up below
#    it) and drop the now-relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$paraCtor = $d.Paragraphs.Item(16)
$ctorXml = '<w:p xmlns:w="' + $wmain + '">' +
    '<w:pPr><w:pStyle w:val="Subtitle"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:i/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">public </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>CreateGameView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>)</w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$paraCtor.Range.InsertXML($ctorXml)

# Add a new "public Game createGame( )" subtitle line right after it.
$paraCtor = $d.Paragraphs.Item(16)
$paraCtor.Range.InsertParagraphAfter()
$paraCreateGameSig = $d.Paragraphs.Item(17)
$createGameSigXml = '<w:p xmlns:w="' + $wmain + '">' +
    '<w:pPr><w:pStyle w:val="Subtitle"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:i/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">public </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Game </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>createGame</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>)</w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$paraCreateGameSig.Range.InsertXML($createGameSigXml)

# ---------------------------------------------------------------------------
# 3) "Method Writeups" heading: collapse the two runs (with the spellStart/
#    spellEnd proofing marks around "Writeups") into one clean run.
# ---------------------------------------------------------------------------
$paraHeading = $d.Paragraphs.Item(19)
$headingXml = '<w:p xmlns:w="' + $wmain + '">' +
    '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
    '<w:r><w:t>Method Writeups</w:t></w:r>' +
    '</w:p>'
$paraHeading.Range.InsertXML($headingXml)

# ---------------------------------------------------------------------------
# 4) Replace the placeholder "This method does some stuff with two strings"
#    paragraph with the real CreateGameView() write-up, then add the
#    createGame() method write-up (Heading4 signature + description) after it.
# ---------------------------------------------------------------------------
$paraBody = $d.Paragraphs.Item(21)
$bodyXml = '<w:p xmlns:w="' + $wmain + '"><w:r><w:t>Instantiate this view.</w:t></w:r></w:p>'
$paraBody.Range.InsertXML($bodyXml)

$paraBody = $d.Paragraphs.Item(21)
$paraBody.Range.InsertParagraphAfter()
$paraCreateGameHeading = $d.Paragraphs.Item(22)
$createGameHeadingXml = '<w:p xmlns:w="' + $wmain + '">' +
    '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">public </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Game </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>createGame</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '</w:p>'
$paraCreateGameHeading.Range.InsertXML($createGameHeadingXml)

$paraCreateGameHeading = $d.Paragraphs.Item(22)
$paraCreateGameHeading.Range.InsertParagraphAfter()
$paraCreateGameBody = $d.Paragraphs.Item(23)
$createGameBodyXml = '<w:p xmlns:w="' + $wmain + '"><w:r><w:t>Compiles the settings in the view and creates a Game object.</w:t></w:r></w:p>'
$paraCreateGameBody.Range.InsertXML($createGameBodyXml)
